# Apply the row-data rotation/swap described by the diff.
# Columns A, B, D, E, F, G, H, Q, R, AC, AX hold the per-observation data
# that gets shuffled between rows while all other columns (shared
# metadata for the locality) stay put.
#
# Row groups affected:
#   104 -> 105 -> 107 -> 106 -> 104   (4-cycle)
#   115 <-> 116                       (swap)
#   122 <-> 123                       (swap)
#   129 <-> 130                       (swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","Q","R","AC","AX")

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$row").Value()
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $data[$c]
    }
}

# Snapshot the "before" data for every affected row first, since writes
# must not influence subsequent reads.
$d104 = Get-RowData 104
$d105 = Get-RowData 105
$d106 = Get-RowData 106
$d107 = Get-RowData 107

$d115 = Get-RowData 115
$d116 = Get-RowData 116

$d122 = Get-RowData 122
$d123 = Get-RowData 123

$d129 = Get-RowData 129
$d130 = Get-RowData 130

# Apply the 4-cycle: new104 = old105, new105 = old107, new107 = old106, new106 = old104
Set-RowData 104 $d105
Set-RowData 105 $d107
Set-RowData 107 $d106
Set-RowData 106 $d104

# Apply the simple swaps
Set-RowData 115 $d116
Set-RowData 116 $d115

Set-RowData 122 $d123
Set-RowData 123 $d122

Set-RowData 129 $d130
Set-RowData 130 $d129
